$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E8").Value = 1
$ws.Range("E11").Value = 1
$ws.Range("E13").Value = 2

$ws.Range("E13").Select()
